# Append a new data row (row 2) to Sheet1:
#   A2 = "vickey"  (text)
#   B2 = "1"       (text, not the number 1)
#
# B2 must be stored as text, so force the cell's number format to
# "Text" (@) before assigning the value - otherwise Excel's COM layer
# auto-coerces the numeric-looking string "1" into a real number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "vickey"

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "1"
